$wb = $excel.ActiveWorkbook

# Sheet: Single| alpha = 0
$ws = $wb.Worksheets.Item("Single| alpha = 0")
$c = $ws.Range("B3")
$c.NumberFormat = "@"
$c.Value = "28.74"
$c.ClearFormats()
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "13.53"
$c.ClearFormats()
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = "10.32"
$c.ClearFormats()
$c = $ws.Range("F3")
$c.NumberFormat = "@"
$c.Value = "9.08"
$c.ClearFormats()
$c = $ws.Range("H3")
$c.NumberFormat = "@"
$c.Value = "5.23"
$c.ClearFormats()
$c = $ws.Range("B4")
$c.NumberFormat = "@"
$c.Value = "28.74"
$c.ClearFormats()
$c = $ws.Range("C4")
$c.NumberFormat = "@"
$c.Value = "52.31"
$c.ClearFormats()
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "65.84"
$c.ClearFormats()
$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = "76.15"
$c.ClearFormats()
$c = $ws.Range("F4")
$c.NumberFormat = "@"
$c.Value = "85.23"
$c.ClearFormats()
$c = $ws.Range("G4")
$c.NumberFormat = "@"
$c.Value = "91.83"
$c.ClearFormats()
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = "0.74"
$c.ClearFormats()
$c = $ws.Range("F7")
$c.NumberFormat = "@"
$c.Value = "0.64"
$c.ClearFormats()
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.36"
$c.ClearFormats()
$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = "-0.20"
$c.ClearFormats()
$c = $ws.Range("F9")
$c.NumberFormat = "@"
$c.Value = "-0.26"
$c.ClearFormats()
$c = $ws.Range("G9")
$c.NumberFormat = "@"
$c.Value = "-0.30"
$c.ClearFormats()
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = "-0.13"
$c.ClearFormats()

# Sheet: Single| alpha = 0.2
$ws = $wb.Worksheets.Item("Single| alpha = 0.2")
$c = $ws.Range("B3")
$c.NumberFormat = "@"
$c.Value = "29.20"
$c.ClearFormats()
$c = $ws.Range("C3")
$c.NumberFormat = "@"
$c.Value = "23.64"
$c.ClearFormats()
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "13.34"
$c.ClearFormats()
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = "10.45"
$c.ClearFormats()
$c = $ws.Range("F3")
$c.NumberFormat = "@"
$c.Value = "8.79"
$c.ClearFormats()
$c = $ws.Range("B4")
$c.NumberFormat = "@"
$c.Value = "29.20"
$c.ClearFormats()
$c = $ws.Range("C4")
$c.NumberFormat = "@"
$c.Value = "52.85"
$c.ClearFormats()
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "66.19"
$c.ClearFormats()
$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = "76.63"
$c.ClearFormats()
$c = $ws.Range("F4")
$c.NumberFormat = "@"
$c.Value = "85.42"
$c.ClearFormats()
$c = $ws.Range("F8")
$c.NumberFormat = "@"
$c.Value = "-0.44"
$c.ClearFormats()
$c = $ws.Range("H10")
$c.NumberFormat = "@"
$c.Value = "0.08"
$c.ClearFormats()
$c = $ws.Range("H12")
$c.NumberFormat = "@"
$c.Value = "-0.26"
$c.ClearFormats()

# Sheet: Single| alpha = 0.5
$ws = $wb.Worksheets.Item("Single| alpha = 0.5")
$c = $ws.Range("B3")
$c.NumberFormat = "@"
$c.Value = "29.36"
$c.ClearFormats()
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = "10.43"
$c.ClearFormats()
$c = $ws.Range("H3")
$c.NumberFormat = "@"
$c.Value = "5.17"
$c.ClearFormats()
$c = $ws.Range("B4")
$c.NumberFormat = "@"
$c.Value = "29.36"
$c.ClearFormats()
$c = $ws.Range("C4")
$c.NumberFormat = "@"
$c.Value = "53.10"
$c.ClearFormats()
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "66.36"
$c.ClearFormats()
$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = "76.78"
$c.ClearFormats()
$c = $ws.Range("F4")
$c.NumberFormat = "@"
$c.Value = "85.52"
$c.ClearFormats()
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = "0.29"
$c.ClearFormats()
$c = $ws.Range("H5")
$c.NumberFormat = "@"
$c.Value = "-0.37"
$c.ClearFormats()
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.29"
$c.ClearFormats()
$c = $ws.Range("H8")
$c.NumberFormat = "@"
$c.Value = "0.20"
$c.ClearFormats()

# Sheet: Single| alpha = 0.7
$ws = $wb.Worksheets.Item("Single| alpha = 0.7")
$c = $ws.Range("B3")
$c.NumberFormat = "@"
$c.Value = "29.39"
$c.ClearFormats()
$c = $ws.Range("G3")
$c.NumberFormat = "@"
$c.Value = "6.41"
$c.ClearFormats()
$c = $ws.Range("B4")
$c.NumberFormat = "@"
$c.Value = "29.39"
$c.ClearFormats()
$c = $ws.Range("C4")
$c.NumberFormat = "@"
$c.Value = "53.17"
$c.ClearFormats()
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "66.39"
$c.ClearFormats()
$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = "76.79"
$c.ClearFormats()
$c = $ws.Range("F4")
$c.NumberFormat = "@"
$c.Value = "85.53"
$c.ClearFormats()
$c = $ws.Range("G5")
$c.NumberFormat = "@"
$c.Value = "0.59"
$c.ClearFormats()
$c = $ws.Range("H8")
$c.NumberFormat = "@"
$c.Value = "0.20"
$c.ClearFormats()
$c = $ws.Range("F11")
$c.NumberFormat = "@"
$c.Value = "0.05"
$c.ClearFormats()
$c = $ws.Range("G12")
$c.NumberFormat = "@"
$c.Value = "-0.27"
$c.ClearFormats()

# Sheet: Single| alpha = 0.9
$ws = $wb.Worksheets.Item("Single| alpha = 0.9")
$c = $ws.Range("B2")
$c.NumberFormat = "@"
$c.Value = "2.37"
$c.ClearFormats()
$c = $ws.Range("B3")
$c.NumberFormat = "@"
$c.Value = "29.40"
$c.ClearFormats()
$c = $ws.Range("G3")
$c.NumberFormat = "@"
$c.Value = "6.39"
$c.ClearFormats()
$c = $ws.Range("B4")
$c.NumberFormat = "@"
$c.Value = "29.40"
$c.ClearFormats()
$c = $ws.Range("C4")
$c.NumberFormat = "@"
$c.Value = "53.21"
$c.ClearFormats()
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "66.40"
$c.ClearFormats()
$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = "76.78"
$c.ClearFormats()
$c = $ws.Range("G4")
$c.NumberFormat = "@"
$c.Value = "91.92"
$c.ClearFormats()
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = "0.30"
$c.ClearFormats()
$c = $ws.Range("F6")
$c.NumberFormat = "@"
$c.Value = "-0.25"
$c.ClearFormats()
$c = $ws.Range("G7")
$c.NumberFormat = "@"
$c.Value = "-0.00"
$c.ClearFormats()
$c = $ws.Range("H8")
$c.NumberFormat = "@"
$c.Value = "0.20"
$c.ClearFormats()
$c = $ws.Range("G12")
$c.NumberFormat = "@"
$c.Value = "-0.27"
$c.ClearFormats()

# Sheet: Single| alpha = 1.0
$ws = $wb.Worksheets.Item("Single| alpha = 1.0")
$c = $ws.Range("C4")
$c.NumberFormat = "@"
$c.Value = "53.23"
$c.ClearFormats()
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "66.41"
$c.ClearFormats()
$c = $ws.Range("H10")
$c.NumberFormat = "@"
$c.Value = "0.07"
$c.ClearFormats()
$c = $ws.Range("G12")
$c.NumberFormat = "@"
$c.Value = "-0.27"
$c.ClearFormats()

# Sheet: Pair| alpha = 0.2
$ws = $wb.Worksheets.Item("Pair| alpha = 0.2")
$c = $ws.Range("B3")
$c.NumberFormat = "@"
$c.Value = "48.34"
$c.ClearFormats()
$c = $ws.Range("B4")
$c.NumberFormat = "@"
$c.Value = "48.34"
$c.ClearFormats()
$c = $ws.Range("C5")
$c.NumberFormat = "@"
$c.Value = "-0.32"
$c.ClearFormats()
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "-0.74"
$c.ClearFormats()
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.09"
$c.ClearFormats()

# Sheet: Pair| alpha = 0.9
$ws = $wb.Worksheets.Item("Pair| alpha = 0.9")
$c = $ws.Range("C5")
$c.NumberFormat = "@"
$c.Value = "-0.15"
$c.ClearFormats()

# Sheet: Pair| alpha = 1.0
$ws = $wb.Worksheets.Item("Pair| alpha = 1.0")
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "25.10"
$c.ClearFormats()
$c = $ws.Range("C4")
$c.NumberFormat = "@"
$c.Value = "74.90"
$c.ClearFormats()
